$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in match scores (Home score in column C, Away score in column E)
# for the 3rd-8th October 2021 matchdays (rows 10-21), using the "different
# prediction method" values from the commit.
$scores = @{
    10 = @(1, 1)
    11 = @(2, 2)
    12 = @(1, 2)
    13 = @(0, 0)
    14 = @(2, 2)
    15 = @(3, 0)
    16 = @(1, 1)
    17 = @(2, 0)
    18 = @(4, 2)
    19 = @(1, 0)
    20 = @(2, 2)
    21 = @(2, 1)
}

foreach ($row in ($scores.Keys | Sort-Object)) {
    $pair = $scores[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}

# Update the view state: scroll/selection moved to G20
$ws.Activate()
$ws.Range("G20").Select()
